$d = $word.ActiveDocument

# Helper to run a simple Find & ReplaceAll over the whole document content.
function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-All "2025-12-20 Saturday" "2025-12-21 Sunday"

# Table cell values that are unique across the document -- safe to
# Find/ReplaceAll.
Replace-All "652×6=3912" "836×2=1672"
Replace-All "554×7=3878" "913×9=8217"
Replace-All "974×7=6818" "952×4=3808"
Replace-All "207×6=1242" "751×4=3004"
Replace-All "220×2=440" "900×3=2700"
Replace-All "470×9=4230" "150×7=1050"
Replace-All "906×2=1812" "701×7=4907"
Replace-All "361×4=1444" "113×6=678"
Replace-All "794×5=3970" "399×6=2394"
Replace-All "533×8=4264" "588×5=2940"
Replace-All "604×7=4228" "160×2=320"
Replace-All "610×6=3660" "424×7=2968"
Replace-All "928×4=3712" "239×5=1195"
Replace-All "169×7=1183" "386×4=1544"
Replace-All "396×3=1188" "910×2=1820"
Replace-All "407×4=1628" "580×6=3480"
Replace-All "401×8=3208" "306×2=612"
Replace-All "335×6=2010" "747×4=2988"
Replace-All "315×8=2520" "748×8=5984"
Replace-All "744×3=2232" "872×7=6104"
Replace-All "253×4=1012" "744×9=6696"
Replace-All "388×3=1164" "310×3=930"
Replace-All "398×7=2786" "350×8=2800"

# "158×9=1422" appears twice in the source table (row 10 col 2, and
# row 20 col 4) and each occurrence becomes a *different* value, so a
# blanket text replace can't be used here -- address the two cells
# directly through the table/row/column indices instead.
$t = $d.Tables.Item(1)
$t.Cell(10, 2).Range.Text = "673×4=2692"
$t.Cell(20, 4).Range.Text = "674×9=6066"

Write-Host "edit complete"
